$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column width adjustments (C: 15 -> ~13.140625, E: 21.42578125 -> 41)
$ws.Columns.Item(3).ColumnWidth = 12.333333333333334
$ws.Columns.Item(5).ColumnWidth = 40.166666666666664

# 1) Write all NEW shared-string cells FIRST, in the exact order needed to reproduce
#    the target sharedStrings.xml table order (indices 16..38).
$ws.Cells.Item(6,1).Value2 = 'Sample Hack Data'  # A6
$ws.Cells.Item(6,10).Value2 = ' 0.068482 s'  # J6
$ws.Cells.Item(6,12).Value2 = '0.001121 s'  # L6
$ws.Cells.Item(3,10).Value2 = '0.069484 s'  # J3
$ws.Cells.Item(3,12).Value2 = '0.000955 s'  # L3
$ws.Cells.Item(4,12).Value2 = '0.000945 s'  # L4
$ws.Cells.Item(5,10).Value2 = '0.065585 s'  # J5
$ws.Cells.Item(5,12).Value2 = '0.000213 s'  # L5
$ws.Cells.Item(4,10).Value2 = '0.075676 s'  # J4
$ws.Cells.Item(3,1).Value2 = 'Credit Card Customer Data'  # A3
$ws.Cells.Item(4,1).Value2 = 'student_performance'  # A4
$ws.Cells.Item(5,1).Value2 = 'Mall_Customers'  # A5
$ws.Cells.Item(2,11).Value2 = ' 22.145ms'  # K2
$ws.Cells.Item(2,13).Value2 = '3.3849ms'  # M2
$ws.Cells.Item(3,11).Value2 = '49.3692ms'  # K3
$ws.Cells.Item(3,13).Value2 = ' 3.5976ms'  # M3
$ws.Cells.Item(4,5).Value2 = ' skipped (dataset too large: 14003 samples)'  # E4
$ws.Cells.Item(4,11).Value2 = '672.6516ms'  # K4
$ws.Cells.Item(4,13).Value2 = '73.075ms'  # M4
$ws.Cells.Item(5,11).Value2 = ' 5.5642ms'  # K5
$ws.Cells.Item(5,13).Value2 = '894.5µs'  # M5
$ws.Cells.Item(6,11).Value2 = '33.0354ms'  # K6
$ws.Cells.Item(6,13).Value2 = ' 6.7235ms'  # M6

# 2) Write all NEW numeric-only cells (order irrelevant; sheetData is always
#    re-serialised sorted by row/column).
$ws.Cells.Item(2,3).Value2 = 36528388099.321999
$ws.Cells.Item(2,5).Value2 = 0.72563148587909398
$ws.Cells.Item(2,7).Value2 = 0.519030794865175
$ws.Cells.Item(2,9).Value2 = 365.569495897016
$ws.Cells.Item(3,2).Value2 = 748453645657.82495
$ws.Cells.Item(3,3).Value2 = 760529014052.61694
$ws.Cells.Item(3,4).Value2 = 0.62962938004549396
$ws.Cells.Item(3,5).Value2 = 0.50080512590602899
$ws.Cells.Item(3,6).Value2 = 0.62556195507549295
$ws.Cells.Item(3,7).Value2 = 0.92850522883973197
$ws.Cells.Item(3,8).Value2 = 542.71569921305399
$ws.Cells.Item(3,9).Value2 = 523.65122680305899
$ws.Cells.Item(4,2).Value2 = 7166053.9505460002
$ws.Cells.Item(4,3).Value2 = 7166058.8818704505
$ws.Cells.Item(4,4).Value2 = 0.27627451802386799
$ws.Cells.Item(4,6).Value2 = 1.42521565904897
$ws.Cells.Item(4,7).Value2 = 1.42521377459166
$ws.Cells.Item(4,8).Value2 = 6464.3486941232604
$ws.Cells.Item(4,9).Value2 = 6464.3346108914302
$ws.Cells.Item(5,2).Value2 = 387114.74027399998
$ws.Cells.Item(5,3).Value2 = 387114.74027402699
$ws.Cells.Item(5,4).Value2 = 0.47926366436618301
$ws.Cells.Item(5,5).Value2 = 0.47926366436618401
$ws.Cells.Item(5,6).Value2 = 0.76591411115234498
$ws.Cells.Item(5,7).Value2 = 0.76591411115234498
$ws.Cells.Item(5,8).Value2 = 300.97646672732401
$ws.Cells.Item(5,9).Value2 = 300.97646672732401
$ws.Cells.Item(6,2).Value2 = 6915190.6272700001
$ws.Cells.Item(6,3).Value2 = 6915190.6272704499
$ws.Cells.Item(6,4).Value2 = 0.64161761953179597
$ws.Cells.Item(6,5).Value2 = 0.64161761953179497
$ws.Cells.Item(6,6).Value2 = 0.46420287440166702
$ws.Cells.Item(6,7).Value2 = 0.46420287440162999
$ws.Cells.Item(6,8).Value2 = 985.08343154885495
$ws.Cells.Item(6,9).Value2 = 985.08343154885597

# Update selection to F31
$ws.Range("F31").Select() | Out-Null
